$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1345.2391
$ws.Range("I15").Value = 1345.2391
$ws.Range("K15").Value = 4035.7173
$ws.Range("M15").Value = -3866.7173

$ws.Range("H69").Value = 3500
$ws.Range("I69").Value = 3500
$ws.Range("K69").Value = 10500
$ws.Range("M69").Value = -9626

$ws.Range("H72").Value = 3500
$ws.Range("I72").Value = 3500
$ws.Range("K72").Value = 31500
$ws.Range("M72").Value = -27132

$ws.Range("H86").Value = 4546
$ws.Range("I86").Value = 4495.5
$ws.Range("J86").Value = 4579.6665
$ws.Range("K86").Value = 4495.5
$ws.Range("L86").Value = 4579.6665
$ws.Range("M86").Value = -3372.5
$ws.Range("N86").Value = -6825.6665

$ws.Range("H89").Value = 4546
$ws.Range("I89").Value = 4495.5
$ws.Range("J89").Value = 4579.6665
$ws.Range("K89").Value = 22477.5
$ws.Range("L89").Value = 22898.3325
$ws.Range("M89").Value = -16861.5
$ws.Range("N89").Value = -34130.3325

$ws.Range("H110").Value = 42000
$ws.Range("J110").Value = 42000
$ws.Range("L110").Value = 42000
$ws.Range("N110").Value = -50180

$ws.Range("H137").Value = 2284.0264
$ws.Range("I137").Value = 1302.08
$ws.Range("K137").Value = 3906.24
$ws.Range("M137").Value = -1356.24

$ws.Range("H138").Value = 4129.516
$ws.Range("J138").Value = 4259.963
$ws.Range("L138").Value = 12779.889
$ws.Range("N138").Value = -23059.889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1304
$ws.Range("I2").Value = 1007.6
$ws.Range("K2").Value = 1007.6
$ws.Range("M2").Value = -894.6

$ws.Range("H74").Value = 1715.2
$ws.Range("I74").Value = 1281.5
$ws.Range("J74").Value = 3450
$ws.Range("K74").Value = 1281.5
$ws.Range("L74").Value = 3450
$ws.Range("M74").Value = -407.5
$ws.Range("N74").Value = -5198

$ws.Range("H77").Value = 1715.2
$ws.Range("I77").Value = 1281.5
$ws.Range("J77").Value = 3450
$ws.Range("K77").Value = 6407.5
$ws.Range("L77").Value = 17250
$ws.Range("M77").Value = -2039.5
$ws.Range("N77").Value = -25986

$ws.Range("H116").Value = 1304
$ws.Range("I116").Value = 1007.6
$ws.Range("K116").Value = 1007.6
$ws.Range("M116").Value = 1286.4

$ws.Range("H122").Value = 1803.1818
$ws.Range("I122").Value = 1493.5
$ws.Range("K122").Value = 4480.5
$ws.Range("M122").Value = -2030.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1304
$ws.Range("I3").Value = 1007.6
$ws.Range("K3").Value = 1007.6
$ws.Range("M3").Value = -893.6

$ws.Range("H20").Value = 10000
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3454.389
$ws.Range("I31").Value = 2318.25
$ws.Range("K31").Value = 2318.25
$ws.Range("M31").Value = -2023.25

$ws.Range("H34").Value = 3454.389
$ws.Range("I34").Value = 2318.25
$ws.Range("K34").Value = 2318.25
$ws.Range("M34").Value = -2116.25

$ws.Range("H52").Value = 131100
$ws.Range("J52").Value = 134800
$ws.Range("L52").Value = 134800
$ws.Range("N52").Value = -135388

$ws.Range("H55").Value = 6000
$ws.Range("J55").Value = 6000
$ws.Range("L55").Value = 6000
$ws.Range("N55").Value = -6630

$ws.Range("H86").Value = 20687.35
$ws.Range("I86").Value = 9645.6
$ws.Range("J86").Value = 31729.1
$ws.Range("K86").Value = 9645.6
$ws.Range("L86").Value = 31729.1
$ws.Range("M86").Value = -8522.6
$ws.Range("N86").Value = -33975.1

$ws.Range("H89").Value = 20687.35
$ws.Range("I89").Value = 9645.6
$ws.Range("J89").Value = 31729.1
$ws.Range("K89").Value = 48228
$ws.Range("L89").Value = 158645.5
$ws.Range("M89").Value = -42612
$ws.Range("N89").Value = -169877.5

$ws.Range("H107").Value = 2362.3845
$ws.Range("I107").Value = 1341.25
$ws.Range("K107").Value = 1341.25
$ws.Range("M107").Value = 578.75

$ws.Range("H122").Value = 2101.8333
$ws.Range("I122").Value = 652.75
$ws.Range("K122").Value = 1958.25
$ws.Range("M122").Value = 491.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 215.27272
$ws.Range("I12").Value = 264.66666
$ws.Range("K12").Value = 793.9999799999999
$ws.Range("M12").Value = -620.9999799999999

$ws.Range("H38").Value = 93.8
$ws.Range("I38").Value = 62.5
$ws.Range("J38").Value = 114.666664
$ws.Range("K38").Value = 187.5
$ws.Range("L38").Value = 343.999992
$ws.Range("M38").Value = 159.5
$ws.Range("N38").Value = -1037.999992

$ws.Range("H68").Value = 2980
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 2980
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H137").Value = 4651.2144
$ws.Range("I137").Value = 3799.6667
$ws.Range("J137").Value = 4883.4546
$ws.Range("K137").Value = 11399.0001
$ws.Range("L137").Value = 14650.3638
$ws.Range("M137").Value = -6299.000100000001
$ws.Range("N137").Value = -24850.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 5000
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -4730
$ws.Range("N70").Value = -5540

$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 5000
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -4064
$ws.Range("N73").Value = -6872

$ws.Range("H80").Value = 5186.143
$ws.Range("I80").Value = 4382.8335
$ws.Range("J80").Value = 10006
$ws.Range("K80").Value = 4382.8335
$ws.Range("L80").Value = 10006
$ws.Range("M80").Value = -3384.8335
$ws.Range("N80").Value = -12002

$ws.Range("H83").Value = 5186.143
$ws.Range("I83").Value = 4382.8335
$ws.Range("J83").Value = 10006
$ws.Range("K83").Value = 21914.1675
$ws.Range("L83").Value = 50030
$ws.Range("M83").Value = -16922.1675
$ws.Range("N83").Value = -60014

$ws.Range("H97").Value = 395.55173
$ws.Range("I97").Value = 430.78262
$ws.Range("K97").Value = 430.78262
$ws.Range("M97").Value = 65.21737999999999

$ws.Range("H122").Value = 1718.3334
$ws.Range("I122").Value = 1718.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5155.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2705.0002
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2400
$ws.Range("J46").Value = 2400
$ws.Range("L46").Value = 2400
$ws.Range("N46").Value = -2776

$ws.Range("H55").Value = 387
$ws.Range("I55").Value = 387
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 387
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -214
$ws.Range("N55").ClearContents()

$ws.Range("H61").Value = 4575
$ws.Range("I61").Value = 4600
$ws.Range("J61").Value = 4500
$ws.Range("K61").Value = 4600
$ws.Range("L61").Value = 4500
$ws.Range("M61").Value = -4398
$ws.Range("N61").Value = -4904

$ws.Range("H82").Value = 2474
$ws.Range("I82").Value = 2474
$ws.Range("K82").Value = 2474
$ws.Range("M82").Value = -2113

$ws.Range("H85").Value = 2474
$ws.Range("I85").Value = 2474
$ws.Range("K85").Value = 2474
$ws.Range("M85").Value = -1226

$ws.Range("H93").Value = 4117.1665
$ws.Range("I93").Value = 4117.1665
$ws.Range("K93").Value = 4117.1665
$ws.Range("M93").Value = -2869.1665

$ws.Range("H100").Value = 1122.25
$ws.Range("I100").Value = 496.33334
$ws.Range("K100").Value = 496.33334
$ws.Range("M100").Value = 44.66665999999998

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H113").Value = 4575
$ws.Range("I113").Value = 4600
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 4600
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -2430
$ws.Range("N113").Value = -8840

$ws.Range("H132").Value = 4956.8335
$ws.Range("I132").Value = 4609.5557
$ws.Range("K132").Value = 13828.6671
$ws.Range("M132").Value = -11298.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7079.8335
$ws.Range("I81").Value = 3870.5
$ws.Range("J81").Value = 13498.5
$ws.Range("K81").Value = 7741
$ws.Range("L81").Value = 26997
$ws.Range("M81").Value = -6680
$ws.Range("N81").Value = -29119

$ws.Range("H84").Value = 7079.8335
$ws.Range("I84").Value = 3870.5
$ws.Range("J84").Value = 13498.5
$ws.Range("K84").Value = 38705
$ws.Range("L84").Value = 134985
$ws.Range("M84").Value = -33401
$ws.Range("N84").Value = -145593

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
